# Logged Week 16 and performed season sim from Week 17
# Update Road (row 3) Target Depth Data totals on both the OFF and DEF sheets

$wb = $excel.ActiveWorkbook

# --- Offense sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 356
$wsOff.Range("C3").Value = 237
$wsOff.Range("D3").Value = 62
$wsOff.Range("E3").Value = 22

# --- Defense sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 434
$wsDef.Range("C3").Value = 311
$wsDef.Range("D3").Value = 94
$wsDef.Range("E3").Value = 48
